# Commit: "update table from cruise to month"
# The model changed from using factor "Cruise" (level "OR1-1242") to using
# factor "Month" (level "October"); factor order is now Depth, DRM, Month.
# This renames / reorders the coefficient rows & columns across the
# coefficients, sw, coefmat.full and coefmat.subset sheets, and updates the
# model id on the msTable sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "coefficients": row 1 holds term names across columns C:H, rows 2-3
# hold the "full"/"subset" coefficient estimates for those terms.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("coefficients")

$ws.Cells.Item(1, 3).Value = "Depth"
$ws.Cells.Item(1, 4).Value = "DRM"
$ws.Cells.Item(1, 5).Value = "MonthOctober"
$ws.Cells.Item(1, 6).Value = "Depth:DRM"
$ws.Cells.Item(1, 7).Value = "Depth:MonthOctober"
$ws.Cells.Item(1, 8).Value = "DRM:MonthOctober"

# row 2 ("full")
$ws.Cells.Item(2, 3).Value = -0.06320796689958592
$ws.Cells.Item(2, 4).Value = -0.1512151265170543
$ws.Cells.Item(2, 5).Value = -0.2412501141555634
$ws.Cells.Item(2, 6).Value = 0.04916820716799926
$ws.Cells.Item(2, 7).Value = 0.1237872147162948
$ws.Cells.Item(2, 8).Value = 0.1554659802726933

# row 3 ("subset")
$ws.Cells.Item(3, 3).Value = -0.06320796689958594
$ws.Cells.Item(3, 4).Value = -0.1512151265170543
$ws.Cells.Item(3, 5).Value = -0.2412501141555634
$ws.Cells.Item(3, 6).Value = 0.06179045099957722
$ws.Cells.Item(3, 7).Value = 0.1237872147162948
$ws.Cells.Item(3, 8).Value = 0.1554659802726933

# ---------------------------------------------------------------------------
# Sheet "msTable": the winning model's id changed (A3). It must stay text
# (like "123456" in A2), so force text format while assigning, then restore
# the original (default) cell style so no real formatting change sticks.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("msTable")
$origStyle = $ws.Cells.Item(3, 1).Style
$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "12356"
$ws.Cells.Item(3, 1).Style = $origStyle

# ---------------------------------------------------------------------------
# Sheet "sw": term names in column A (A2:A6); A7 ("Depth:DRM") is unchanged.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("sw")
$ws.Cells.Item(2, 1).Value = "Depth"
$ws.Cells.Item(3, 1).Value = "DRM"
$ws.Cells.Item(4, 1).Value = "Month"
$ws.Cells.Item(5, 1).Value = "Depth:Month"
$ws.Cells.Item(6, 1).Value = "DRM:Month"

# ---------------------------------------------------------------------------
# Sheet "coefmat.full": rows 3-8 hold one term per row (name + stats).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("coefmat.full")

$ws.Cells.Item(3, 1).Value = "Depth"
$ws.Cells.Item(3, 2).Value = -0.06320796689958594
$ws.Cells.Item(3, 3).Value = 0.0315480445101673
$ws.Cells.Item(3, 4).Value = 0.0329548908406748
$ws.Cells.Item(3, 5).Value = 1.918014755538837
$ws.Cells.Item(3, 6).Value = 0.05510914

$ws.Cells.Item(4, 1).Value = "DRM"
$ws.Cells.Item(4, 2).Value = -0.1512151265170543
$ws.Cells.Item(4, 3).Value = 0.03669034352477185
$ws.Cells.Item(4, 4).Value = 0.03816622791508788
$ws.Cells.Item(4, 5).Value = 3.962013926382176
$ws.Cells.Item(4, 6).Value = 0.00007432000000000001

$ws.Cells.Item(5, 1).Value = "MonthOctober"
$ws.Cells.Item(5, 2).Value = -0.2412501141555634
$ws.Cells.Item(5, 3).Value = 0.03875367753173927
$ws.Cells.Item(5, 4).Value = 0.04033787180746528
$ws.Cells.Item(5, 5).Value = 5.980734811867679
$ws.Cells.Item(5, 6).Value = 0

$ws.Cells.Item(6, 1).Value = "Depth:DRM"
$ws.Cells.Item(6, 2).Value = 0.04916820716799927
$ws.Cells.Item(6, 3).Value = 0.0343806908504884
$ws.Cells.Item(6, 4).Value = 0.03518685583538469
$ws.Cells.Item(6, 5).Value = 1.397345855453064
$ws.Cells.Item(6, 6).Value = 0.16230959

$ws.Cells.Item(7, 1).Value = "Depth:MonthOctober"
$ws.Cells.Item(7, 2).Value = 0.1237872147162948
$ws.Cells.Item(7, 3).Value = 0.04070155971880979
$ws.Cells.Item(7, 4).Value = 0.04249420047344463
$ws.Cells.Item(7, 5).Value = 2.913037857804893
$ws.Cells.Item(7, 6).Value = 0.00357931

$ws.Cells.Item(8, 1).Value = "DRM:MonthOctober"
$ws.Cells.Item(8, 2).Value = 0.1554659802726933
$ws.Cells.Item(8, 3).Value = 0.03832298143311283
$ws.Cells.Item(8, 4).Value = 0.04012850263089879
$ws.Cells.Item(8, 5).Value = 3.874203373663513
$ws.Cells.Item(8, 6).Value = 0.00010697

# ---------------------------------------------------------------------------
# Sheet "coefmat.subset": same term renames, its own coefficient values.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("coefmat.subset")

$ws.Cells.Item(3, 1).Value = "Depth"
$ws.Cells.Item(3, 2).Value = -0.06320796689958594
$ws.Cells.Item(3, 3).Value = 0.0315480445101673
$ws.Cells.Item(3, 4).Value = 0.0329548908406748
$ws.Cells.Item(3, 5).Value = 1.918014755538837
$ws.Cells.Item(3, 6).Value = 0.05510914

$ws.Cells.Item(4, 1).Value = "DRM"
$ws.Cells.Item(4, 2).Value = -0.1512151265170543
$ws.Cells.Item(4, 3).Value = 0.03669034352477185
$ws.Cells.Item(4, 4).Value = 0.03816622791508788
$ws.Cells.Item(4, 5).Value = 3.962013926382176
$ws.Cells.Item(4, 6).Value = 0.00007432000000000001

$ws.Cells.Item(5, 1).Value = "MonthOctober"
$ws.Cells.Item(5, 2).Value = -0.2412501141555634
$ws.Cells.Item(5, 3).Value = 0.03875367753173927
$ws.Cells.Item(5, 4).Value = 0.04033787180746528
$ws.Cells.Item(5, 5).Value = 5.980734811867679
$ws.Cells.Item(5, 6).Value = 0

$ws.Cells.Item(6, 1).Value = "Depth:DRM"
$ws.Cells.Item(6, 2).Value = 0.06179045099957722
$ws.Cells.Item(6, 3).Value = 0.02656207345534019
$ws.Cells.Item(6, 4).Value = 0.02785720776989634
$ws.Cells.Item(6, 5).Value = 2.21811358517958
$ws.Cells.Item(6, 6).Value = 0.02654709

$ws.Cells.Item(7, 1).Value = "Depth:MonthOctober"
$ws.Cells.Item(7, 2).Value = 0.1237872147162948
$ws.Cells.Item(7, 3).Value = 0.04070155971880979
$ws.Cells.Item(7, 4).Value = 0.04249420047344463
$ws.Cells.Item(7, 5).Value = 2.913037857804893
$ws.Cells.Item(7, 6).Value = 0.00357931

$ws.Cells.Item(8, 1).Value = "DRM:MonthOctober"
$ws.Cells.Item(8, 2).Value = 0.1554659802726933
$ws.Cells.Item(8, 3).Value = 0.03832298143311283
$ws.Cells.Item(8, 4).Value = 0.04012850263089879
$ws.Cells.Item(8, 5).Value = 3.874203373663513
$ws.Cells.Item(8, 6).Value = 0.00010697
